# Apply "Add data for 2022-06-07" update:
#  - rename sheet/title to reflect new "through" date (05-29 -> 05-30)
#  - update the shared string used for the May row label
#  - update May row (row 6) and Total row (row 7) values for columns B..I

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab (also updates <sheet name="..."> in workbook.xml)
$ws.Name = "Through 2022-05-30"

# Update the "May (through 05-29)" label -> "May (through 05-30)"
$ws.Range("A6").Value = "May (through 05-30)"

# Updated counts for row 6 (May) and row 7 (Total), columns B..I
$mayValues = @(19, 43, 59, 48, 44, 69, 105, 110)
$totalValues = @(108, 205, 312, 294, 199, 331, 628, 661)

for ($i = 0; $i -lt $mayValues.Length; $i++) {
    $col = 2 + $i  # column B = 2
    $ws.Cells.Item(6, $col).Value = $mayValues[$i]
    $ws.Cells.Item(7, $col).Value = $totalValues[$i]
}
